$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Diagonal cells: keep their original text, but each one now lives one row
# lower than before (row N -> row N+1), and a new top row (row 2) plus the
# remainder of each row is filled with "-".
$diag = @{
    3 = @{ col = "A"; text = "Bane=Disguise Self" }
    4 = @{ col = "B"; text = "Mirror Image=Ray of Enfeeblement" }
    5 = @{ col = "C"; text = "Bestow Curse=Vampiric Touch" }
    6 = @{ col = "D"; text = "Death Ward=Dimension Door" }
    7 = @{ col = "E"; text = "Dominate Person=Modify Memory" }
}

$cols = @("A", "B", "C", "D", "E")

# Copy the existing style (font/format) from A1 so every newly created cell
# reuses the same style index instead of minting a new one.
$ws.Range("A1").Copy()

for ($row = 2; $row -le 7; $row++) {
    if ($row -eq 2) {
        $startCol = 0
    } else {
        $startCol = $row - 3
    }
    for ($ci = $startCol; $ci -le 4; $ci++) {
        $col = $cols[$ci]
        $cellRef = "$col$row"
        $ws.Range($cellRef).PasteSpecial(-4122)
        if ($diag.ContainsKey($row) -and $diag[$row].col -eq $col) {
            $ws.Range($cellRef).Value = $diag[$row].text
        } else {
            $ws.Range($cellRef).Value = "-"
        }
    }
}
